$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------------
# 1. "Attempt Delivery" (row 12) now also covers the old "Update Delivery
#    Status" step, so its description gains the DELIVERED/UNDELIVERED
#    outcomes and the row is marked Complete instead of In Progress.
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "When a shipment is marked as RECEIVED_AT_DEST_BRANCH, an agent from the destination branch attempts the delivery and performs below task:" + $nl + "1. Updates shipment table such that status = OUT_FOR_DELIVERY or DELIVERED or UNDELIVERED" + $nl + "2. Insert a record into shipment_tracker; shipment_id, agent_id, creation_datetime,status=OUT_FOR_DELIVERY or DELIVERED or UNDELIVERED" + $nl

$ws.Range("F12").Value = "Complete"

# Re-colour row 12 from the "In Progress" yellow fill to the "Complete"
# green fill used elsewhere in the Core section (reuse existing formats
# instead of fabricating new ones).
$ws.Range("B10").Copy()
$ws.Range("B12").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C12:E12").PasteSpecial(-4122)

$ws.Range("F10").Copy()
$ws.Range("F12").PasteSpecial(-4122)

$ws.Range("G10").Copy()
$ws.Range("G12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. The old "Update Delivery Status" row (13) is now folded into row 12
#    above, so remove it entirely - this shifts "Track Shipment" up to
#    row 13 and shrinks the A8:A14 merge down to A8:A13.
# ---------------------------------------------------------------------------
$ws.Rows(13).Delete()

# ---------------------------------------------------------------------------
# 3. Update the saved view/selection to match the now-shorter sheet.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("B12:G12").Select() | Out-Null
